$d = $word.ActiveDocument

# wdParagraph unit for Range.Expand
$wdParagraph = 4

# --- Locate the block of reference hyperlinks that sits right before the
# "References:" heading (the 9 plain hyperlink paragraphs that start with
# the vNegKfkopoQ YouTube link and end with the first mongodb-and-django
# hyperlink). ---
$rngBlock1Start = $d.Content.Duplicate
$foundStart1 = $rngBlock1Start.Find.Execute("https://www.youtube.com/watch?v=vNegKfkopoQ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngBlock1Start.Expand($wdParagraph) | Out-Null

$rngBlock1End = $d.Content.Duplicate
$foundEnd1 = $rngBlock1End.Find.Execute("https://www.mongodb.com/compatibility/mongodb-and-django", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngBlock1End.Expand($wdParagraph) | Out-Null

# --- Locate the second block: everything from the "References:" heading
# through to the very end of the document body (before the section
# properties). ---
$rngBlock2Start = $d.Content.Duplicate
$foundStart2 = $rngBlock2Start.Find.Execute("References:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngBlock2Start.Expand($wdParagraph) | Out-Null

$docEnd = $d.Content.End

# Delete the later block first so the earlier block's character offsets
# remain valid.
if ($foundStart2) {
    $rng2 = $d.Range($rngBlock2Start.Start, $docEnd)
    $rng2.Delete()
}

if ($foundStart1 -and $foundEnd1) {
    $rng1 = $d.Range($rngBlock1Start.Start, $rngBlock1End.End)
    $rng1.Delete()
}
